$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain text (dotted/grouped numbers, not real
# numeric values). Force text format before writing so numeric-looking
# strings (e.g. "85.97") are not auto-converted to numbers by Excel, then
# restore the default "Normal" style so the saved cell style matches the
# original (unstyled) cells.
$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range("D2").Value = '27.393.20'
$ws.Range("E2").Value = '  -3.72%  '
$ws.Range("D3").Value = '1.854.31'
$ws.Range("E3").Value = '  -5.01%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.59%  '
$ws.Range("D5").Value = '320.13'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  -0.67%  '
$ws.Range("D7").Value = '0.4474'
$ws.Range("E7").Value = '  -6.32%  '
$ws.Range("D8").Value = '0.3842'
$ws.Range("E8").Value = '  -4.56%  '
$ws.Range("D9").Value = '46.93'
$ws.Range("E9").Value = '  -12.46%  '
$ws.Range("D10").Value = '0.07966'
$ws.Range("E10").Value = '  -6.10%  '
$ws.Range("D11").Value = '1.017'
$ws.Range("E11").Value = '  -3.90%  '
$ws.Range("D12").Value = '21.31'
$ws.Range("E12").Value = '  -3.46%  '
$ws.Range("D13").Value = '1.886.12'
$ws.Range("E13").Value = '  -3.01%  '
$ws.Range("D14").Value = '5.860'
$ws.Range("E14").Value = '  -5.37%  '
$ws.Range("D15").Value = '7.111'
$ws.Range("E15").Value = '  -6.76%  '
$ws.Range("D16").Value = '1.008'
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '85.97'
$ws.Range("E17").Value = '  -3.48%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.00001032'
$ws.Range("E18").Value = '  -4.28%  '
$ws.Range("D19").Value = '0.06525'
$ws.Range("E19").Value = '  -1.17%  '
$ws.Range("D20").Value = '17.08'
$ws.Range("E20").Value = '  -8.88%  '
$ws.Range("D21").Value = '1.005'
$ws.Range("E21").Value = '  -0.78%  '
$ws.Range("D22").Value = '5.498'
$ws.Range("E22").Value = '  -5.49%  '
$ws.Range("D23").Value = '27.410.01'
$ws.Range("E23").Value = '  -3.77%  '
$ws.Range("D24").Value = '10.83'
$ws.Range("E24").Value = '  -5.99%  '
$ws.Range("D25").Value = '2.283'
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("D26").Value = '2.101.12'
$ws.Range("E26").Value = '  -3.72%  '
$ws.Range("D27").Value = '151.06'
$ws.Range("E27").Value = '  -1.92%  '
$ws.Range("D28").Value = '19.37'
$ws.Range("E28").Value = '  -3.99%  '
$ws.Range("D29").Value = '5.525'
$ws.Range("E29").Value = '  -7.23%  '
$ws.Range("D30").Value = '2.029'
$ws.Range("E30").Value = '  -5.95%  '
$ws.Range("D31").Value = '120.47'
$ws.Range("E31").Value = '  -2.51%  '
$ws.Range("D32").Value = '0.09381'
$ws.Range("E32").Value = '  -1.93%  '
$ws.Range("D33").Value = '1.478'
$ws.Range("E33").Value = '  +2.16%  '
$ws.Range("D34").Value = '0.9270'
$ws.Range("E34").Value = '  -6.80%  '
$ws.Range("D35").Value = '3.623'
$ws.Range("E35").Value = '  -1.04%  '
$ws.Range("D36").Value = '5.271'
$ws.Range("E36").Value = '  -5.83%  '
$ws.Range("D37").Value = '0.02228'
$ws.Range("E37").Value = '  -4.82%  '
$ws.Range("D38").Value = '1.223'
$ws.Range("E38").Value = '  -3.09%  '
$ws.Range("D39").Value = '0.05959'
$ws.Range("E39").Value = '  -4.31%  '
$ws.Range("D40").Value = '8.323'
$ws.Range("E40").Value = '  -5.31%  '
$ws.Range("D41").Value = '1.004'
$ws.Range("E41").Value = '  -0.73%  '
$ws.Range("D42").Value = '0.5909'
$ws.Range("E42").Value = '  -5.33%  '
$ws.Range("D43").Value = '0.1853'
$ws.Range("E43").Value = '  -3.93%  '
$ws.Range("D44").Value = '10.25'
$ws.Range("E44").Value = '  -7.74%  '
$ws.Range("E45").Value = '  -3.91%  '
$ws.Range("D46").Value = '0.5646'
$ws.Range("E46").Value = '  -5.47%  '
$ws.Range("D47").Value = '12.24'
$ws.Range("E47").Value = '  -5.36%  '
$ws.Range("D48").Value = '1.923'
$ws.Range("E48").Value = '  -6.87%  '
$ws.Range("D49").Value = '3.348'
$ws.Range("E49").Value = '  -1.61%  '
$ws.Range("D50").Value = '0.06847'
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("E51").Value = '  -0.91%  '

$ws.Range("D2:D50").Style = "Normal"
